$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J -- copy formatting from an existing
# header cell (H1) so the new headers pick up the same bold/border/center
# style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2..44 for columns I (col 9) and J (col 10)
$data = @(
    @(7,7),
    @(8,8),
    @(5,5),
    @(1,1),
    @(1,1),
    @(7,7),
    @(8,8),
    @(6,7),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,7),
    @(7,7),
    @(4,4),
    @(1,1),
    @(1,2),
    @(8,8),
    @(1,2),
    @(15,15),
    @(8,8),
    @(8,8),
    @(1,1),
    @(5,6),
    @(10,10),
    @(6,6),
    @(6,6),
    @(9,9),
    @(8,8),
    @(8,8),
    @(7,7),
    @(9,9),
    @(7,8),
    @(9,9),
    @(6,6),
    @(6,7),
    @(5,5),
    @(6,6),
    @(8,8),
    @(6,6),
    @(4,4),
    @(3,3),
    @(7,7),
    @(2,2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
